# Adds a new record for Doncaster's new provider (A9N8Z) into the
# project_name_lookup table, inserted in sorted position as row 14
# (pushing the existing rows 14-32 down to 15-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new blank row at row 14 (within the table body) and grow the
# table definition to match, so the new row becomes part of Table1 and
# the sheet dimension / autofilter extend to G33.
$ws.Rows("14:14").Insert()
$tbl.Resize($ws.Range("A1:G33"))

# Populate the new record.
$ws.Range("A14").Value = "Doncaster"
$ws.Range("B14").Value = "A9N8Z"
$ws.Range("C14").Value = "Doncaster"
$ws.Range("D14").Value = "0"
$ws.Range("F14").Value = "Doncaster"
$ws.Range("G14").Value = "New TLHC provider"

# Leave the selection where the editor last left it.
$ws.Range("F15").Select() | Out-Null
